$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row labels (column A): EUR, USD, USDEUR
$ws.Range("A2").Value = "EUR"
$ws.Range("A3").Value = "USD"
$ws.Range("A4").Value = "USDEUR"

# Clear old row 5/6 labels and data (CMS spread correlations removed)
$ws.Range("A5").ClearContents()
$ws.Range("A6").ClearContents()

# Header row: formulas pulling labels from column A
$ws.Range("B1").Formula = "=A2"
$ws.Range("C1").Formula = "=A3"
$ws.Range("D1").Formula = "=A4"
$ws.Range("E1").ClearContents()
$ws.Range("F1").ClearContents()

# New correlation values
$ws.Range("B3").Value = 0.7196
$ws.Range("B4").Value = -0.0645
$ws.Range("C4").Value = 0.0763

# Mirror formulas for the symmetric matrix
$ws.Range("C2").Formula = "=B3"
$ws.Range("D2").Formula = "=B4"
$ws.Range("D3").Formula = "=C4"

# Clear now-unused cells (previously held CMS spread correlations)
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("B5:F6").ClearContents()

$ws.Range("C5").Select()
